$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 18.50795745849609
$ws.Range("D2").Value = 259

$ws.Range("C3").Value = 16.96395874023438
$ws.Range("D3").Value = 174

$ws.Range("C4").Value = 16.97492599487305
$ws.Range("D4").Value = 176

$ws.Range("C5").Value = 17.08006858825684
$ws.Range("D5").Value = 174

$ws.Range("C6").Value = 16.94798469543457
$ws.Range("D6").Value = 174
